$wb = $excel.ActiveWorkbook

# --- Sheet N1_D40 ---
$ws = $wb.Worksheets.Item("N1_D40")
$ws.Range("E2").Value = 0.508
$ws.Range("F2").Value = 16.77
$ws.Range("F3").Value = 15.39
$ws.Range("E4").Value = 0.019
$ws.Range("E5").Value = 0.019
$ws.Range("F5").Value = 15.26
$ws.Range("E6").Value = 0.019
$ws.Range("F6").Value = 15.39
$ws.Range("F7").Value = 15.45
$ws.Range("E8").Value = 0.019
$ws.Range("F8").Value = 15.36
$ws.Range("F9").Value = 15.4
$ws.Range("F10").Value = 15.36
$ws.Range("F11").Value = 15.37
$ws.Range("E12").Value = 0.0674
$ws.Range("F12").Value = 15.507

# --- Sheet N1_D60 ---
$ws = $wb.Worksheets.Item("N1_D60")
$ws.Range("E2").Value = 0.028
$ws.Range("F2").Value = 14.86
$ws.Range("F3").Value = 14.92
$ws.Range("F4").Value = 14.86
$ws.Range("E5").Value = 0.028
$ws.Range("F5").Value = 14.83
$ws.Range("F6").Value = 14.78
$ws.Range("F7").Value = 14.79
$ws.Range("E8").Value = 0.028
$ws.Range("F8").Value = 14.85
$ws.Range("E9").Value = 0.027
$ws.Range("F9").Value = 14.9
$ws.Range("E10").Value = 0.026
$ws.Range("F10").Value = 14.86
$ws.Range("F11").Value = 14.9
$ws.Range("E12").Value = 0.0269
$ws.Range("F12").Value = 14.855

# --- Sheet N1_D80 ---
$ws = $wb.Worksheets.Item("N1_D80")
$ws.Range("F2").Value = 18.35
$ws.Range("E3").Value = 0.041
$ws.Range("F3").Value = 18.39
$ws.Range("E4").Value = 0.043
$ws.Range("F4").Value = 18.42
$ws.Range("F5").Value = 18.27
$ws.Range("F6").Value = 18.39
$ws.Range("F7").Value = 18.28
$ws.Range("F8").Value = 18.29
$ws.Range("E9").Value = 0.039
$ws.Range("F9").Value = 18.28
$ws.Range("E10").Value = 0.04
$ws.Range("F10").Value = 18.27
$ws.Range("E11").Value = 0.041
$ws.Range("F11").Value = 18.39
$ws.Range("E12").Value = 0.04039999999999999
$ws.Range("F12").Value = 18.33300000000001

# --- Sheet N1_D100 ---
$ws = $wb.Worksheets.Item("N1_D100")
$ws.Range("E3").Value = 0.061
$ws.Range("F3").Value = 16.4
$ws.Range("E4").Value = 0.058
$ws.Range("F4").Value = 16.51
$ws.Range("E5").Value = 0.059
$ws.Range("F5").Value = 16.44
$ws.Range("E6").Value = 0.059
$ws.Range("F6").Value = 16.32
$ws.Range("E7").Value = 0.059
$ws.Range("F7").Value = 16.33
$ws.Range("E8").Value = 0.059
$ws.Range("F8").Value = 16.32
$ws.Range("F9").Value = 16.32
$ws.Range("F10").Value = 16.44
$ws.Range("E11").Value = 0.059
$ws.Range("F11").Value = 16.44
$ws.Range("E12").Value = 0.05899999999999998
$ws.Range("F12").Value = 16.399
